$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1 (18:00 -> 19:17)
$ws.Range("A1").Value = "Datos actualizados a 15 de Julio de 2020 a las 19:17"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 3572740
$ws.Range("C4").Value = 27663
$ws.Range("D4").Value = 1608996
$ws.Range("E4").Value = 1824223
$ws.Range("G4").Value = 378
$ws.Range("H4").Value = 139521

# Row 5: Brasil
$ws.Range("B5").Value = 1939167
$ws.Range("C5").Value = 7963
$ws.Range("E5").Value = 651210
$ws.Range("G5").Value = 183
$ws.Range("H5").Value = 74445

# Row 6: India
$ws.Range("B6").Value = 968117
$ws.Range("C6").Value = 30630
$ws.Range("D6").Value = 612782
$ws.Range("E6").Value = 330420
$ws.Range("G6").Value = 600
$ws.Range("H6").Value = 24915

# Row 9: Chile
$ws.Range("B9").Value = 321205
$ws.Range("C9").Value = 1712
$ws.Range("D9").Value = 292085
$ws.Range("E9").Value = 21934
$ws.Range("G9").Value = 117
$ws.Range("H9").Value = 7186

# Row 11: España
$ws.Range("B11").Value = 304574
$ws.Range("C11").Value = 875
$ws.Range("G11").Value = 4
$ws.Range("H11").Value = 28413

# Row 13: Reino Unido
$ws.Range("B13").Value = 291911
$ws.Range("C13").Value = 538
$ws.Range("G13").Value = 85
$ws.Range("H13").Value = 45053

# Row 16: Italia
$ws.Range("B16").Value = 243506
$ws.Range("C16").Value = 162
$ws.Range("D16").Value = 196016
$ws.Range("E16").Value = 12493
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = 34997

# Row 19: Alemania
$ws.Range("B19").Value = 200895
$ws.Range("C19").Value = 129
$ws.Range("E19").Value = 5750
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = 9145

# Row 43: Portugal -> Republica Dominicana
$ws.Range("A43").Value = "Republica Dominicana"
$ws.Range("B43").Value = 47671
$ws.Range("C43").Value = 1366
$ws.Range("D43").Value = 23459
$ws.Range("E43").Value = 23283
$ws.Range("G43").Value = 19
$ws.Range("H43").Value = 929

# Row 44: Singapur -> Portugal
$ws.Range("A44").Value = "Portugal"
$ws.Range("B44").Value = 47426
$ws.Range("C44").Value = 375
$ws.Range("D44").Value = 32110
$ws.Range("E44").Value = 13640
$ws.Range("G44").Value = 8
$ws.Range("H44").Value = 1676

# Row 45: Republica Dominicana -> Singapur
$ws.Range("A45").Value = "Singapur"
$ws.Range("B45").Value = 46878
$ws.Range("C45").Value = 249
$ws.Range("D45").Value = 42988
$ws.Range("E45").Value = 3863
$ws.Range("H45").Value = 27

# Row 46: Israel
$ws.Range("B46").Value = 43668
$ws.Range("C46").Value = 1308
$ws.Range("D46").Value = 19894
$ws.Range("E46").Value = 23399

# Row 56: Azerbaiyan -> Irlanda
$ws.Range("A56").Value = "Irlanda"
$ws.Range("B56").Value = 25683
$ws.Range("C56").Value = 13
$ws.Range("D56").Value = 23364
$ws.Range("E56").Value = 571
$ws.Range("G56").Value = 2
$ws.Range("H56").Value = 1748

# Row 57: Irlanda -> Azerbaiyan
$ws.Range("A57").Value = "Azerbaiyan"
$ws.Range("B57").Value = 25672
$ws.Range("C57").Value = 559
$ws.Range("D57").Value = 16695
$ws.Range("E57").Value = 8651
$ws.Range("G57").Value = 7
$ws.Range("H57").Value = 326

# Row 61: Moldavia
$ws.Range("B61").Value = 20040
$ws.Range("C61").Value = 332
$ws.Range("E61").Value = 6083
$ws.Range("G61").Value = 4
$ws.Range("H61").Value = 659

# Row 69: Chequia
$ws.Range("B69").Value = 13421
$ws.Range("C69").Value = 80
$ws.Range("D69").Value = 8505
$ws.Range("E69").Value = 4561

# Row 70: Dinamarca
$ws.Range("B70").Value = 13092
$ws.Range("C70").Value = 31
$ws.Range("D70").Value = 12182
$ws.Range("E70").Value = 300

# Row 80: Costa Rica -> Republica de Macedonia
$ws.Range("A80").Value = "Republica de Macedonia"
$ws.Range("B80").Value = 8530
$ws.Range("C80").Value = 198
$ws.Range("D80").Value = 4565
$ws.Range("E80").Value = 3572
$ws.Range("G80").Value = 4
$ws.Range("H80").Value = 393

# Row 81: Senegal -> Costa Rica
$ws.Range("A81").Value = "Costa Rica"
$ws.Range("B81").Value = 8482
$ws.Range("C81").Value = 0
$ws.Range("D81").Value = 2441
$ws.Range("E81").Value = 6004
$ws.Range("G81").Value = 1
$ws.Range("H81").Value = 37

# Row 82: Republica de Macedonia -> Senegal
$ws.Range("A82").Value = "Senegal"
$ws.Range("B82").Value = 8369
$ws.Range("C82").Value = 126
$ws.Range("D82").Value = 5605
$ws.Range("E82").Value = 2611
$ws.Range("G82").Value = 3
$ws.Range("H82").Value = 153

# Row 96: Luxemburgo
$ws.Range("B96").Value = 5122
$ws.Range("C96").Value = 66
$ws.Range("D96").Value = 4247
$ws.Range("E96").Value = 764

# Row 101: Grecia
$ws.Range("B101").Value = 3910
$ws.Range("C101").Value = 27
$ws.Range("E101").Value = 2343

# Row 111: Malaui -> Libano
$ws.Range("A111").Value = "Libano"
$ws.Range("B111").Value = 2542
$ws.Range("C111").Value = 91
$ws.Range("D111").Value = 1455
$ws.Range("E111").Value = 1049
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = 38

# Row 112: Libano -> Malaui
$ws.Range("A112").Value = "Malaui"
$ws.Range("B112").Value = 2497
$ws.Range("D112").Value = 795
$ws.Range("E112").Value = 1662
$ws.Range("H112").Value = 40

# Row 113: Cuba -> Mali
$ws.Range("A113").Value = "Mali"
$ws.Range("B113").Value = 2433
$ws.Range("C113").Value = 10
$ws.Range("D113").Value = 1764
$ws.Range("E113").Value = 548
$ws.Range("H113").Value = 121

# Row 114: Mali -> Cuba
$ws.Range("A114").Value = "Cuba"
$ws.Range("B114").Value = 2432
$ws.Range("D114").Value = 2275
$ws.Range("E114").Value = 70
$ws.Range("H114").Value = 87

# Row 133: Tunez -> Mozambique
$ws.Range("A133").Value = "Mozambique"
$ws.Range("B133").Value = 1330
$ws.Range("C133").Value = 62
$ws.Range("D133").Value = 375
$ws.Range("E133").Value = 946
$ws.Range("H133").Value = 9

# Row 134: Montenegro -> Tunez
$ws.Range("A134").Value = "Tunez"
$ws.Range("B134").Value = 1319
$ws.Range("C134").Value = 13
$ws.Range("D134").Value = 1091
$ws.Range("E134").Value = 178
$ws.Range("H134").Value = 50

# Row 135: Mozambique -> Montenegro
$ws.Range("A135").Value = "Montenegro"
$ws.Range("B135").Value = 1287
$ws.Range("D135").Value = 330
$ws.Range("E135").Value = 933
$ws.Range("H135").Value = 24

# Row 136: Jordania
$ws.Range("B136").Value = 1201
$ws.Range("C136").Value = 3
$ws.Range("D136").Value = 1016

# Row 164: Birmania
$ws.Range("D164").Value = 266
$ws.Range("E164").Value = 65

# Row 209: Islas Malvinas -> Groenlandia
$ws.Range("A209").Value = "Groenlandia"

# Row 210: Groenlandia -> Islas Malvinas
$ws.Range("A210").Value = "Islas Malvinas"
